$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Plain text / non-numeric-looking values: set directly.
$ws.Range('D2').Value = '30.036.64'
$ws.Range('E2').Value = '  -2.12%  '
$ws.Range('D3').Value = '1.829.85'
$ws.Range('E4').Value = '  +0.00%  '
$ws.Range('E5').Value = '  -4.29%  '
$ws.Range('E7').Value = '  -4.56%  '
$ws.Range('E8').Value = '  -6.46%  '
$ws.Range('E9').Value = '  -5.79%  '
$ws.Range('D10').Value = '1.830.82'
$ws.Range('E10').Value = '  -3.18%  '
$ws.Range('E11').Value = '  -1.79%  '
$ws.Range('E12').Value = '  -5.15%  '
$ws.Range('E13').Value = '  -4.78%  '
$ws.Range('E14').Value = '  -6.75%  '
$ws.Range('E15').Value = '  -7.99%  '
$ws.Range('D16').Value = '29.981.92'
$ws.Range('E16').Value = '  -2.25%  '
$ws.Range('E17').Value = '  +0.05%  '
$ws.Range('E18').Value = '  -3.76%  '
$ws.Range('E19').Value = '  -0.04%  '
$ws.Range('B20').Value = 'ShibaInu'
$ws.Range('C20').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('E20').Value = '  -5.23%  '
$ws.Range('B21').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C21').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D21').Value = '2.074.30'
$ws.Range('E21').Value = '  -5.66%  '
$ws.Range('E22').Value = '  -8.41%  '
$ws.Range('E23').Value = '  -9.36%  '
$ws.Range('E24').Value = '  -6.67%  '
$ws.Range('E25').Value = '  -3.49%  '
$ws.Range('E26').Value = '  -3.27%  '
$ws.Range('E27').Value = '  -6.99%  '
$ws.Range('E28').Value = '  -7.07%  '
$ws.Range('E29').Value = '  -1.03%  '
$ws.Range('E30').Value = '  -2.38%  '
$ws.Range('E31').Value = '  -6.94%  '
$ws.Range('E32').Value = '  -7.60%  '
$ws.Range('E33').Value = '  -6.61%  '
$ws.Range('E34').Value = '  -7.55%  '
$ws.Range('E35').Value = '  -9.25%  '
$ws.Range('E36').Value = '  +0.06%  '
$ws.Range('E37').Value = '  -0.85%  '
$ws.Range('E38').Value = '  -4.82%  '
$ws.Range('E39').Value = '  -1.80%  '
$ws.Range('E40').Value = '  -4.64%  '
$ws.Range('B41').Value = 'PaxDollar'
$ws.Range('C41').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('E41').Value = '  -0.90%  '
$ws.Range('B42').Value = 'RenderToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('E42').Value = '  -8.67%  '
$ws.Range('E43').Value = '  -4.40%  '
$ws.Range('E44').Value = '  -4.77%  '
$ws.Range('E45').Value = '  -8.15%  '
$ws.Range('E47').Value = '  -7.97%  '
$ws.Range('E48').Value = '  -8.84%  '
$ws.Range('E49').Value = '  -7.16%  '
$ws.Range('E50').Value = '  -2.67%  '
$ws.Range('E51').Value = '  -4.99%  '

# Numeric-looking values that must remain TEXT (matches source inlineStr cells):
# Temporarily force a text number-format so Excel does not coerce the literal
# into a real number, then restore the cell's original style so no stray
# formatting change is left behind.
$textCells = @{
    'D4' = '1.001'
    'D5' = '228.68'
    'D6' = '1.001'
    'D7' = '0.4601'
    'D8' = '0.2679'
    'D9' = '0.06179'
    'D11' = '0.07342'
    'D12' = '15.85'
    'D13' = '4.864'
    'D15' = '0.6143'
    'D18' = '223.61'
    'D19' = '1.001'
    'D20' = '0.000007191'
    'D22' = '12.20'
    'D23' = '4.790'
    'D24' = '5.810'
    'D25' = '164.28'
    'D26' = '9.035'
    'D27' = '17.48'
    'D28' = '1.828'
    'D29' = '0.1004'
    'D30' = '1.367'
    'D31' = '4.048'
    'D32' = '3.728'
    'D33' = '0.04764'
    'D34' = '1.124'
    'D35' = '0.6903'
    'D37' = '2.690'
    'D38' = '0.01787'
    'D39' = '2.604'
    'D40' = '0.8788'
    'D41' = '0.9957'
    'D42' = '1.896'
    'D43' = '102.46'
    'D44' = '5.412'
    'D45' = '0.3959'
    'D46' = '6.842'
    'D48' = '58.74'
    'D49' = '8.355'
    'D50' = '0.05520'
    'D51' = '32.27'
}
foreach ($addr in $textCells.Keys) {
    $cell = $ws.Range($addr)
    $savedStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $textCells[$addr]
    $cell.Style = $savedStyle
}
